# Deploy the implementation guide:
#  - refresh the "Date" metadata value
#  - replace the retired "REPORT" data-type concept with "EXOMIZER"
#  - add the new "IGV", "CNVVIZ" and "COVGENE" data-type concepts
#  - keep the existing "OTHER" concept, and append the new "FUSG" concept

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet: bump the publication Date ------------------------
$ws1.Range("B8").Value = "2023-03-09T16:02:20+00:00"

# --- Concepts sheet -----------------------------------------------------
# Row 8 used to hold the "REPORT" / "Data Report" concept; the FHIR
# CodeSystem renamed/replaced it with "EXOMIZER" / "Exomizer Report".
$ws2.Range("B8").Value = "EXOMIZER"
$ws2.Range("C8").Value = "Exomizer Report"

# Three brand-new concepts are inserted right after it. Duplicate the
# formatting of the row above (via Copy, which carries over the border /
# fill / text styling) before overwriting the Code & Display text so the
# new rows look identical to the existing table rows.
$ws2.Range("A8:D8").Copy($ws2.Range("A9:D9"))
$ws2.Range("B9").Value = "IGV"
$ws2.Range("C9").Value = "IGV Track"

$ws2.Range("A8:D8").Copy($ws2.Range("A10:D10"))
$ws2.Range("B10").Value = "CNVVIZ"
$ws2.Range("C10").Value = "CNV Visualization"

$ws2.Range("A8:D8").Copy($ws2.Range("A11:D11"))
$ws2.Range("B11").Value = "COVGENE"
$ws2.Range("C11").Value = "Coverage by Gene Report"

# The old "OTHER" / "Undefined Data Type" row (previously row 9) now
# lands on row 12, after the three new concepts.
$ws2.Range("A8:D8").Copy($ws2.Range("A12:D12"))
$ws2.Range("B12").Value = "OTHER"
$ws2.Range("C12").Value = "Undefined Data Type"

# Finally, append the new "FUSG" / "Gene fusion" concept as the last row.
$ws2.Range("A8:D8").Copy($ws2.Range("A13:D13"))
$ws2.Range("B13").Value = "FUSG"
$ws2.Range("C13").Value = "Gene fusion"
